$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2 currently holds "{{AppUsers.AppUserStoreMappings.CodeDraft}}" (drafted/manual code).
# Replace it with the auto-generated code placeholder used in E2.
$ws.Range("D2").Value = "{{AppUsers.AppUserStoreMappings.Code}}"

# Row 1 header no longer needs the taller custom height - restore default row height.
$ws.Rows(1).AutoFit()

# Move the active selection to D2 (the field that changed).
$ws.Range("D2").Select()
